# aggiornamento fino a 28 luglio
# Appends new daily rows (302-328) to the sheet, continuing the existing
# date series in column A (dates 2021-06-29 .. 2021-07-25, serials
# 44376..44402) with zero values in columns B, C and D - matching the
# formatting/style already used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 302
$lastNewRow  = 328
$firstNewSerial = 44376

# Copy the formatting (style) of the last existing data row (A301) onto
# the new A-column cells so the new dates keep the same date number
# format / border / alignment.
$ws.Range("A301").Copy($ws.Range("A$firstNewRow`:A$lastNewRow"))

for ($row = $firstNewRow; $row -le $lastNewRow; $row++) {
    $serial = $firstNewSerial + ($row - $firstNewRow)
    $ws.Cells.Item($row, 1).Value2 = $serial
}

# Columns B, C, D are all 0 for the new rows, plain (unstyled) numbers
# just like the existing rows.
$ws.Range("B$firstNewRow`:D$lastNewRow").Value2 = 0
